# Replace JS libraries note with CDN guidance, add new "non-explicit page name"
# audit row, and add two new reference hyperlinks (matches commit:
# "replacing JS libraries with CDNs").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D16: broaden the "hidden text" advice to also mention alt attributes ---
$ws.Range("D16").Value = "Ne pas mettre de texte caché avec des mots clés, ou dans les attributs alt des images"

# --- F16: new reference link (alt attribute / hidden text article) ---
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.redacteur.com/blog/seo-balise-alt-images/") | Out-Null
$ws.Range("F16").Value = "https://www.redacteur.com/blog/seo-balise-alt-images/"
$ws.Range("F16").Style = $ws.Range("F15").Style

# --- Row 9: brand-new audit line about non-explicit page naming ---
$ws.Range("A9").Value = "SEO"
$ws.Range("B9").Value = "nom de page non-explicite"
$ws.Range("C9").Value = "La page de contact est nommée page-2"
$ws.Range("D9").Value = "Donner un nom de page adapté pour l'URL"
$ws.Range("E9").Value = "Modifier ""page2"" par ""contact"""

# --- F21: new reference link (CDN article, with #:~:text= fragment) ---
$cdnBase = "https://www.keycdn.com/blog/why-use-a-cdn"
$cdnFragment = ":~:text=Faster%20performance%20and%20lower%20latency,-Of%20course%2C%20the&text=Using%20a%20CDN%20allows%20us,both%20static%20and%20dynamic%20content."
$cdnFull = $cdnBase + "#" + $cdnFragment
$ws.Hyperlinks.Add($ws.Range("F21"), $cdnBase, $cdnFragment) | Out-Null
$ws.Range("F21").Value = $cdnFull
$ws.Range("F21").Style = $ws.Range("F15").Style

# --- Selection moved by the editing author ---
$ws.Range("G10").Select()
